$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.611132179096228)
    3  = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.242251378316819)
    4  = @(0.04271373187048222, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 0, 10.69140994080654)
    5  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 3.754798637575387)
    6  = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.611132179096228)
    7  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    8  = @(0.04271373187048222, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1, 0.9931934699307545)
    9  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 1, 15.88780690183548)
    10 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 6.348428708163715)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
